$d = $word.ActiveDocument

$d.Content.Find.Execute("2023-07-06 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-07-07 Friday", 2) | Out-Null
$d.Content.Find.Execute("99-12=87", $true, $false, $false, $false, $false, $true, 1, $false, "29-11=18", 2) | Out-Null
$d.Content.Find.Execute("17-10=7", $true, $false, $false, $false, $false, $true, 1, $false, "96-24=72", 2) | Out-Null
$d.Content.Find.Execute("68-16=52", $true, $false, $false, $false, $false, $true, 1, $false, "46-24=22", 2) | Out-Null
$d.Content.Find.Execute("35+7=42", $true, $false, $false, $false, $false, $true, 1, $false, "15+33=48", 2) | Out-Null
$d.Content.Find.Execute("27+31=58", $true, $false, $false, $false, $false, $true, 1, $false, "38-29=9", 2) | Out-Null
$d.Content.Find.Execute("66-41=25", $true, $false, $false, $false, $false, $true, 1, $false, "17+19=36", 2) | Out-Null
$d.Content.Find.Execute("94-60=34", $true, $false, $false, $false, $false, $true, 1, $false, "35+46=81", 2) | Out-Null
$d.Content.Find.Execute("61-4=57", $true, $false, $false, $false, $false, $true, 1, $false, "31-6=25", 2) | Out-Null
$d.Content.Find.Execute("64+20=84", $true, $false, $false, $false, $false, $true, 1, $false, "34-16=18", 2) | Out-Null
$d.Content.Find.Execute("45+16=61", $true, $false, $false, $false, $false, $true, 1, $false, "8+2=10", 2) | Out-Null
$d.Content.Find.Execute("8-4=4", $true, $false, $false, $false, $false, $true, 1, $false, "5+76=81", 2) | Out-Null
$d.Content.Find.Execute("40+5=45", $true, $false, $false, $false, $false, $true, 1, $false, "26+6=32", 2) | Out-Null
$d.Content.Find.Execute("68-56=12", $true, $false, $false, $false, $false, $true, 1, $false, "50+22=72", 2) | Out-Null
$d.Content.Find.Execute("66-11=55", $true, $false, $false, $false, $false, $true, 1, $false, "96-78=18", 2) | Out-Null
$d.Content.Find.Execute("50-2=48", $true, $false, $false, $false, $false, $true, 1, $false, "8+28=36", 2) | Out-Null
$d.Content.Find.Execute("9+75=84", $true, $false, $false, $false, $false, $true, 1, $false, "36+62=98", 2) | Out-Null
$d.Content.Find.Execute("73+5=78", $true, $false, $false, $false, $false, $true, 1, $false, "7+89=96", 2) | Out-Null
$d.Content.Find.Execute("14+45=59", $true, $false, $false, $false, $false, $true, 1, $false, "25+37=62", 2) | Out-Null
$d.Content.Find.Execute("68-46=22", $true, $false, $false, $false, $false, $true, 1, $false, "41-33=8", 2) | Out-Null
$d.Content.Find.Execute("82-30=52", $true, $false, $false, $false, $false, $true, 1, $false, "65-9=56", 2) | Out-Null
$d.Content.Find.Execute("34-15=19", $true, $false, $false, $false, $false, $true, 1, $false, "18+36=54", 2) | Out-Null
$d.Content.Find.Execute("79-31=48", $true, $false, $false, $false, $false, $true, 1, $false, "67+12=79", 2) | Out-Null
$d.Content.Find.Execute("9+21=30", $true, $false, $false, $false, $false, $true, 1, $false, "4+9=13", 2) | Out-Null
$d.Content.Find.Execute("6+41=47", $true, $false, $false, $false, $false, $true, 1, $false, "1+72=73", 2) | Out-Null
$d.Content.Find.Execute("49-22=27", $true, $false, $false, $false, $false, $true, 1, $false, "73-51=22", 2) | Out-Null
$d.Content.Find.Execute("56+5=61", $true, $false, $false, $false, $false, $true, 1, $false, "20+23=43", 2) | Out-Null
$d.Content.Find.Execute("55-34=21", $true, $false, $false, $false, $false, $true, 1, $false, "59+14=73", 2) | Out-Null
$d.Content.Find.Execute("86-39=47", $true, $false, $false, $false, $false, $true, 1, $false, "92-77=15", 2) | Out-Null
$d.Content.Find.Execute("41+7=48", $true, $false, $false, $false, $false, $true, 1, $false, "76+1=77", 2) | Out-Null
$d.Content.Find.Execute("62-30=32", $true, $false, $false, $false, $false, $true, 1, $false, "48-14=34", 2) | Out-Null
$d.Content.Find.Execute("77-24=53", $true, $false, $false, $false, $false, $true, 1, $false, "17+54=71", 2) | Out-Null
$d.Content.Find.Execute("7+81=88", $true, $false, $false, $false, $false, $true, 1, $false, "52+47=99", 2) | Out-Null
$d.Content.Find.Execute("44+14=58", $true, $false, $false, $false, $false, $true, 1, $false, "48-38=10", 2) | Out-Null
$d.Content.Find.Execute("45-43=2", $true, $false, $false, $false, $false, $true, 1, $false, "19+50=69", 2) | Out-Null
$d.Content.Find.Execute("92+2=94", $true, $false, $false, $false, $false, $true, 1, $false, "81-70=11", 2) | Out-Null
$d.Content.Find.Execute("90-24=66", $true, $false, $false, $false, $false, $true, 1, $false, "75-61=14", 2) | Out-Null
$d.Content.Find.Execute("39+37=76", $true, $false, $false, $false, $false, $true, 1, $false, "35+49=84", 2) | Out-Null
$d.Content.Find.Execute("43-22=21", $true, $false, $false, $false, $false, $true, 1, $false, "81-59=22", 2) | Out-Null
$d.Content.Find.Execute("1+20=21", $true, $false, $false, $false, $false, $true, 1, $false, "4+27=31", 2) | Out-Null
$d.Content.Find.Execute("79-36=43", $true, $false, $false, $false, $false, $true, 1, $false, "60-55=5", 2) | Out-Null
$d.Content.Find.Execute("65+20=85", $true, $false, $false, $false, $false, $true, 1, $false, "55+28=83", 2) | Out-Null
$d.Content.Find.Execute("28+53=81", $true, $false, $false, $false, $false, $true, 1, $false, "14+60=74", 2) | Out-Null
$d.Content.Find.Execute("36+6=42", $true, $false, $false, $false, $false, $true, 1, $false, "50-5=45", 2) | Out-Null
$d.Content.Find.Execute("79-56=23", $true, $false, $false, $false, $false, $true, 1, $false, "96-58=38", 2) | Out-Null
$d.Content.Find.Execute("32+3=35", $true, $false, $false, $false, $false, $true, 1, $false, "69-40=29", 2) | Out-Null
$d.Content.Find.Execute("15+14=29", $true, $false, $false, $false, $false, $true, 1, $false, "98-43=55", 2) | Out-Null
$d.Content.Find.Execute("36-15=21", $true, $false, $false, $false, $false, $true, 1, $false, "68-64=4", 2) | Out-Null
$d.Content.Find.Execute("21+57=78", $true, $false, $false, $false, $false, $true, 1, $false, "61+27=88", 2) | Out-Null
$d.Content.Find.Execute("21-20=1", $true, $false, $false, $false, $false, $true, 1, $false, "64-18=46", 2) | Out-Null
$d.Content.Find.Execute("21-6=15", $true, $false, $false, $false, $false, $true, 1, $false, "1+5=6", 2) | Out-Null
$d.Content.Find.Execute("57+8=65", $true, $false, $false, $false, $false, $true, 1, $false, "63-27=36", 2) | Out-Null
$d.Content.Find.Execute("13+48=61", $true, $false, $false, $false, $false, $true, 1, $false, "63+5=68", 2) | Out-Null
$d.Content.Find.Execute("97-77=20", $true, $false, $false, $false, $false, $true, 1, $false, "49-46=3", 2) | Out-Null
$d.Content.Find.Execute("59+8=67", $true, $false, $false, $false, $false, $true, 1, $false, "95-34=61", 2) | Out-Null
$d.Content.Find.Execute("92-32=60", $true, $false, $false, $false, $false, $true, 1, $false, "17+3=20", 2) | Out-Null
$d.Content.Find.Execute("27+11=38", $true, $false, $false, $false, $false, $true, 1, $false, "34+62=96", 2) | Out-Null
$d.Content.Find.Execute("18-1=17", $true, $false, $false, $false, $false, $true, 1, $false, "19+24=43", 2) | Out-Null
$d.Content.Find.Execute("43+3=46", $true, $false, $false, $false, $false, $true, 1, $false, "49-18=31", 2) | Out-Null
$d.Content.Find.Execute("97-0=97", $true, $false, $false, $false, $false, $true, 1, $false, "20+65=85", 2) | Out-Null
$d.Content.Find.Execute("82-61=21", $true, $false, $false, $false, $false, $true, 1, $false, "69-59=10", 2) | Out-Null
$d.Content.Find.Execute("59-36=23", $true, $false, $false, $false, $false, $true, 1, $false, "44-15=29", 2) | Out-Null
$d.Content.Find.Execute("29+14=43", $true, $false, $false, $false, $false, $true, 1, $false, "16+35=51", 2) | Out-Null
$d.Content.Find.Execute("12-5=7", $true, $false, $false, $false, $false, $true, 1, $false, "74-11=63", 2) | Out-Null
$d.Content.Find.Execute("31+35=66", $true, $false, $false, $false, $false, $true, 1, $false, "5+4=9", 2) | Out-Null
$d.Content.Find.Execute("69+15=84", $true, $false, $false, $false, $false, $true, 1, $false, "9-8=1", 2) | Out-Null
$d.Content.Find.Execute("49-14=35", $true, $false, $false, $false, $false, $true, 1, $false, "88-78=10", 2) | Out-Null
$d.Content.Find.Execute("13+6=19", $true, $false, $false, $false, $false, $true, 1, $false, "7+36=43", 2) | Out-Null
$d.Content.Find.Execute("35+41=76", $true, $false, $false, $false, $false, $true, 1, $false, "29+24=53", 2) | Out-Null
$d.Content.Find.Execute("31+20=51", $true, $false, $false, $false, $false, $true, 1, $false, "20+27=47", 2) | Out-Null
$d.Content.Find.Execute("13+41=54", $true, $false, $false, $false, $false, $true, 1, $false, "11+10=21", 2) | Out-Null
$d.Content.Find.Execute("23-0=23", $true, $false, $false, $false, $false, $true, 1, $false, "43-20=23", 2) | Out-Null
$d.Content.Find.Execute("33+16=49", $true, $false, $false, $false, $false, $true, 1, $false, "58-38=20", 2) | Out-Null
$d.Content.Find.Execute("65-59=6", $true, $false, $false, $false, $false, $true, 1, $false, "63+13=76", 2) | Out-Null
$d.Content.Find.Execute("12+22=34", $true, $false, $false, $false, $false, $true, 1, $false, "27-20=7", 2) | Out-Null
$d.Content.Find.Execute("65+6=71", $true, $false, $false, $false, $false, $true, 1, $false, "75-25=50", 2) | Out-Null
$d.Content.Find.Execute("20-13=7", $true, $false, $false, $false, $false, $true, 1, $false, "14+4=18", 2) | Out-Null
$d.Content.Find.Execute("90+7=97", $true, $false, $false, $false, $false, $true, 1, $false, "24+17=41", 2) | Out-Null
$d.Content.Find.Execute("61+5=66", $true, $false, $false, $false, $false, $true, 1, $false, "30+7=37", 2) | Out-Null
$d.Content.Find.Execute("0+15=15", $true, $false, $false, $false, $false, $true, 1, $false, "20+35=55", 2) | Out-Null
$d.Content.Find.Execute("1+51=52", $true, $false, $false, $false, $false, $true, 1, $false, "94-83=11", 2) | Out-Null
$d.Content.Find.Execute("99-50=49", $true, $false, $false, $false, $false, $true, 1, $false, "4+40=44", 2) | Out-Null
$d.Content.Find.Execute("20+48=68", $true, $false, $false, $false, $false, $true, 1, $false, "91-3=88", 2) | Out-Null
$d.Content.Find.Execute("97-97=0", $true, $false, $false, $false, $false, $true, 1, $false, "16+75=91", 2) | Out-Null
$d.Content.Find.Execute("40+21=61", $true, $false, $false, $false, $false, $true, 1, $false, "51-19=32", 2) | Out-Null
$d.Content.Find.Execute("37+50=87", $true, $false, $false, $false, $false, $true, 1, $false, "0+13=13", 2) | Out-Null
$d.Content.Find.Execute("40+25=65", $true, $false, $false, $false, $false, $true, 1, $false, "30+56=86", 2) | Out-Null
$d.Content.Find.Execute("41+58=99", $true, $false, $false, $false, $false, $true, 1, $false, "20+59=79", 2) | Out-Null
$d.Content.Find.Execute("61-6=55", $true, $false, $false, $false, $false, $true, 1, $false, "63+36=99", 2) | Out-Null
$d.Content.Find.Execute("92-12=80", $true, $false, $false, $false, $false, $true, 1, $false, "35+32=67", 2) | Out-Null
$d.Content.Find.Execute("70+19=89", $true, $false, $false, $false, $false, $true, 1, $false, "53-19=34", 2) | Out-Null
$d.Content.Find.Execute("40+35=75", $true, $false, $false, $false, $false, $true, 1, $false, "37-29=8", 2) | Out-Null
$d.Content.Find.Execute("77-38=39", $true, $false, $false, $false, $false, $true, 1, $false, "1+49=50", 2) | Out-Null
$d.Content.Find.Execute("33+38=71", $true, $false, $false, $false, $false, $true, 1, $false, "63+1=64", 2) | Out-Null
$d.Content.Find.Execute("64-35=29", $true, $false, $false, $false, $false, $true, 1, $false, "11+42=53", 2) | Out-Null
$d.Content.Find.Execute("57+29=86", $true, $false, $false, $false, $false, $true, 1, $false, "95-85=10", 2) | Out-Null
$d.Content.Find.Execute("84-81=3", $true, $false, $false, $false, $false, $true, 1, $false, "21-17=4", 2) | Out-Null
$d.Content.Find.Execute("99-4=95", $true, $false, $false, $false, $false, $true, 1, $false, "26+1=27", 2) | Out-Null
$d.Content.Find.Execute("35+25=60", $true, $false, $false, $false, $false, $true, 1, $false, "34+59=93", 2) | Out-Null
$d.Content.Find.Execute("36-33=3", $true, $false, $false, $false, $false, $true, 1, $false, "69-27=42", 2) | Out-Null
$d.Content.Find.Execute("44+53=97", $true, $false, $false, $false, $false, $true, 1, $false, "20+22=42", 2) | Out-Null
